$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.46000000000007
$ws.Range("H2").Value = [double]"3.416070845000481e-16"
$ws.Range("K2").Value = 42.02580948201606
$ws.Range("L2").Value = "[33.69431903704325, 50.357299926988865]"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1.50318447288881
$ws.Range("P2").Value = "[1.2893423303021168, 1.7170266154755032]"
$ws.Range("S2").Value = 59.73285847679992
$ws.Range("T2").Value = "[54.40308250487962, 65.06263444872022]"
$ws.Range("W2").Value = 17.08668668668674
$ws.Range("X2").Value = 16.32228228228233
$ws.Range("Y2").Value = 17.85109109109115

# Row 3 updates
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 25.63000000000057
$ws.Range("H3").Value = [double]"3.416070845000481e-16"
$ws.Range("K3").Value = 38.60139191185458
$ws.Range("L3").Value = "[28.151808528410683, 49.05097529529847]"
$ws.Range("M3").Value = [double]"3.450351115930061e-12"
$ws.Range("N3").Value = [double]"3.450351115930061e-12"
$ws.Range("O3").Value = -0.5031579825569237
$ws.Range("P3").Value = "[-0.779894872963232, -0.22642109215061534]"
$ws.Range("Q3").Value = 0.0004055895124976949
$ws.Range("R3").Value = 0.0004055895124976949
$ws.Range("S3").Value = 55.67030290948367
$ws.Range("T3").Value = "[50.0388590661417, 61.301746752825636]"
$ws.Range("W3").Value = 2.052452452452499
$ws.Range("X3").Value = 0.9236036036036244
$ws.Range("Y3").Value = 3.181301301301374

$wb.Save()
